$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 245.663428745944
$ws.Range("R2").Value = 2210.970858713496
$ws.Range("S2").Value = 0.0351295791722707
$ws.Range("T2").Value = 0.03512957917227069

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 213.2051213984853
$ws.Range("R3").Value = 1918.846092586368
$ws.Range("S3").Value = 0.03048807968827691
$ws.Range("T3").Value = 0.0304880796882769

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 298.0983238102614
$ws.Range("R4").Value = 2682.884914292352
$ws.Range("S4").Value = 0.04262770702530407
$ws.Range("T4").Value = 0.04262770702530407

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 1858.099253027218
$ws.Range("R5").Value = 16722.89327724496
$ws.Range("S5").Value = 0.2657059911292734
$ws.Range("T5").Value = 0.2657059911292733

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 1612.59768633203
$ws.Range("R6").Value = 14513.37917698827
$ws.Range("S6").Value = 0.2305995580384364
$ws.Range("T6").Value = 0.2305995580384363

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 2254.695685182068
$ws.Range("R7").Value = 20292.26116663861
$ws.Range("S7").Value = 0.322418810916675
$ws.Range("T7").Value = 0.3224188109166749

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 165.742592060032
$ws.Range("R8").Value = 1491.683328540288
$ws.Range("S8").Value = 0.02370099424123742
$ws.Range("T8").Value = 0.02370099424123741

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 143.8438339863893
$ws.Range("R9").Value = 1294.594505877504
$ws.Range("S9").Value = 0.02056949778916271
$ws.Range("T9").Value = 0.0205694977891627

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 201.1190234105173
$ws.Range("R10").Value = 1810.071210694656
$ws.Range("S10").Value = 0.02875978199936356
$ws.Range("T10").Value = 0.02875978199936355
